$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The rows 3, 4, 6 and 7 had their weekly records rotated/re-ordered.
# Capture the "before" values of the columns that change (D, L, M, N, O, P, R, S)
# for each of the affected rows, then reassign them following the mapping:
#   new row 3 <- old row 6
#   new row 4 <- old row 3
#   new row 6 <- old row 7
#   new row 7 <- old row 4

function Get-RowData($row) {
    $data = @{}
    $data.D = $ws.Cells.Item($row, 4).Value2
    $data.L = $ws.Cells.Item($row, 12).Value2
    $data.M = $ws.Cells.Item($row, 13).Value2
    $data.N = $ws.Cells.Item($row, 14).Value2
    $data.O = $ws.Cells.Item($row, 15).Value2
    $data.P = $ws.Cells.Item($row, 16).Value2
    $data.R = $ws.Cells.Item($row, 18).Value2
    $data.S = $ws.Cells.Item($row, 19).Value2
    return $data
}

function Set-RowData($row, $data) {
    $ws.Cells.Item($row, 4).Value2 = $data.D
    $ws.Cells.Item($row, 12).Value2 = $data.L
    $ws.Cells.Item($row, 13).Value2 = $data.M
    $ws.Cells.Item($row, 14).Value2 = $data.N
    $ws.Cells.Item($row, 15).Value2 = $data.O
    $ws.Cells.Item($row, 16).Value2 = $data.P
    $ws.Cells.Item($row, 18).Value2 = $data.R
    $ws.Cells.Item($row, 19).Value2 = $data.S
}

$row3 = Get-RowData 3
$row4 = Get-RowData 4
$row6 = Get-RowData 6
$row7 = Get-RowData 7

Set-RowData 3 $row6
Set-RowData 4 $row3
Set-RowData 6 $row7
Set-RowData 7 $row4
